# Insert a new data row after the first data row (row 2), shifting all
# existing rows (old rows 3-48) down by one (new rows 4-49), and populate
# the newly inserted row 3 with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3; Excel shifts rows 3..48 down to 4..49
# and copies formatting from the row above (keeps the date style on D).
$ws.Rows("3:3").Insert()

# Populate the new row 3 with the new record's values.
$ws.Cells.Item(3, 1).Value = 5
$ws.Cells.Item(3, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(3, 3).Value = "Maule"
$ws.Cells.Item(3, 4).Value = 44552
$ws.Cells.Item(3, 5).Value = 7
$ws.Cells.Item(3, 6).Value = "Fruta"
$ws.Cells.Item(3, 7).Value = 100101
$ws.Cells.Item(3, 8).Value = "Berries"
$ws.Cells.Item(3, 9).Value = 100101001
$ws.Cells.Item(3, 10).Value = "Arándano (blue)"
$ws.Cells.Item(3, 11).Value = "Sin especificar"
$ws.Cells.Item(3, 12).Value = "Primera"
$ws.Cells.Item(3, 13).Value = 180
$ws.Cells.Item(3, 14).Value = 4000
$ws.Cells.Item(3, 15).Value = 4000
$ws.Cells.Item(3, 16).Value = 4000
$ws.Cells.Item(3, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(3, 18).Value = "Provincia de Linares"
$ws.Cells.Item(3, 19).Value = 2000
$ws.Cells.Item(3, 20).Value = 2
